$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9998340106875503
$ws.Range("C2").Value = 0.2532704065238534
$ws.Range("D2").Value = 0.02824048846233751
$ws.Range("E2").Value = 0.1184694903237427
$ws.Range("F2").Value = 0.8287855757083591
$ws.Range("L2").Value = 0.1886793512133664
$ws.Range("M2").Value = 0.2225080754306177
$ws.Range("N2").Value = 1.376402873543597
$ws.Range("O2").Value = 2.873585931596779
$ws.Range("B3").Value = 0.9154134562347451
$ws.Range("C3").Value = 0.2429535528171414
$ws.Range("D3").Value = 0.0270744962661098
$ws.Range("E3").Value = 0.1194771694375198
$ws.Range("F3").Value = 0.8230670106637845
$ws.Range("L3").Value = 0.1860235826453902
$ws.Range("M3").Value = 0.209200329540252
$ws.Range("N3").Value = 1.390789325767411
$ws.Range("O3").Value = 2.869311570831997
$ws.Range("B4").Value = 0.8637903451067643
$ws.Range("C4").Value = 0.236565522160376
$ws.Range("D4").Value = 0.02635243062194093
$ws.Range("E4").Value = 0.120133443549435
$ws.Range("F4").Value = 0.820112033225648
$ws.Range("L4").Value = 0.1844843663126738
$ws.Range("M4").Value = 0.2011066445347822
$ws.Range("N4").Value = 1.400155429781734
$ws.Range("O4").Value = 2.868562608234498
$ws.Range("B5").Value = 0.8428078309320028
$ws.Range("C5").Value = 0.2339490491903575
$ws.Range("D5").Value = 0.0260566556715176
$ws.Range("E5").Value = 0.1204103349571142
$ws.Range("F5").Value = 0.8190476355548952
$ws.Range("L5").Value = 0.1838801668474517
$ws.Range("M5").Value = 0.1978280302640059
$ws.Range("N5").Value = 1.404106192046054
$ws.Range("O5").Value = 2.868728632849013
$ws.Range("B6").Value = 0.8393270153734989
$ws.Range("C6").Value = 0.2335137867595876
$ws.Range("D6").Value = 0.02600745072112431
$ws.Range("E6").Value = 0.1204568839050761
$ws.Range("F6").Value = 0.8188793327823305
$ws.Range("L6").Value = 0.1837812332829003
$ws.Range("M6").Value = 0.1972848097824595
$ws.Range("N6").Value = 1.40477030729345
$ws.Range("O6").Value = 2.868784651350097
$ws.Range("B7").Value = 0.8635071461042685
$ws.Range("C7").Value = 0.2365302891510765
$ws.Range("D7").Value = 0.02634844785949397
$ws.Range("E7").Value = 0.1201371395088096
$ws.Range("F7").Value = 0.820097112548126
$ws.Range("L7").Value = 0.1844761244960509
$ws.Range("M7").Value = 0.2010623482783416
$ws.Range("N7").Value = 1.400208168545728
$ws.Range("O7").Value = 2.868562939805258
$ws.Range("B8").Value = 0.9706825803161792
$ws.Range("C8").Value = 0.2497243277053229
$ws.Range("D8").Value = 0.02783973895871128
$ws.Range("E8").Value = 0.1188091501259581
$ws.Range("F8").Value = 0.8266983059276853
$ws.Range("L8").Value = 0.1877446940578125
$ws.Range("M8").Value = 0.2179036217923382
$ws.Range("N8").Value = 1.381252765874144
$ws.Range("O8").Value = 2.871722630734695
$ws.Range("B9").Value = 1.182491285369963
$ws.Range("C9").Value = 0.2751687732997539
$ws.Range("D9").Value = 0.03071480266741844
$ws.Range("E9").Value = 0.1165024395907093
$ws.Range("F9").Value = 0.844062936668692
$ws.Range("L9").Value = 0.194878148232867
$ws.Range("M9").Value = 0.251536570215201
$ws.Range("N9").Value = 1.348306900267978
$ws.Range("O9").Value = 2.892822084706154
$ws.Range("B10").Value = 1.339068263590605
$ws.Range("C10").Value = 0.2935962198797597
$ws.Range("D10").Value = 0.03279637627289844
$ws.Range("E10").Value = 0.1149883120733781
$ws.Range("F10").Value = 0.8595264885622527
$ws.Range("L10").Value = 0.2005588575999013
$ws.Range("M10").Value = 0.2766111450380961
$ws.Range("N10").Value = 1.326675052089094
$ws.Range("O10").Value = 2.917446675793911
$ws.Range("B11").Value = 1.41050093536569
$ws.Range("C11").Value = 0.301920488724619
$ws.Range("D11").Value = 0.03373653787130593
$ws.Range("E11").Value = 0.1143385579396892
$ws.Range("F11").Value = 0.8671514069051227
$ws.Range("L11").Value = 0.2032383865866905
$ws.Range("M11").Value = 0.2880963012296078
$ws.Range("N11").Value = 1.317392492080899
$ws.Range("O11").Value = 2.930638455276096
$ws.Range("B12").Value = 1.437579164946499
$ws.Range("C12").Value = 0.3050641451535796
$ws.Range("D12").Value = 0.03409156606640096
$ws.Range("E12").Value = 0.1140981144268727
$ws.Range("F12").Value = 0.8701238267604623
$ws.Range("L12").Value = 0.2042667260072903
$ws.Range("M12").Value = 0.2924565858272317
$ws.Range("N12").Value = 1.313957642468516
$ws.Range("O12").Value = 2.935920541193383
$ws.Range("B13").Value = 1.431746148260913
$ws.Range("C13").Value = 0.3043874858347806
$ws.Range("D13").Value = 0.0340151487519762
$ws.Range("E13").Value = 0.1141496491175662
$ws.Range("F13").Value = 0.869479879389047
$ws.Range("L13").Value = 0.2040446479903437
$ws.Range("M13").Value = 0.2915170297711924
$ws.Range("N13").Value = 1.314693828891187
$ws.Range("O13").Value = 2.934770192977822
$ws.Range("B14").Value = 1.412728119076348
$ws.Range("C14").Value = 0.3021792917476489
$ws.Range("D14").Value = 0.03376576621605665
$ws.Range("E14").Value = 0.1143186642358565
$ws.Range("F14").Value = 0.8673942448907184
$ws.Range("L14").Value = 0.2033227152142985
$ws.Range("M14").Value = 0.2884548030377232
$ws.Range("N14").Value = 1.31710829631718
$ws.Range("O14").Value = 2.931067268030972
$ws.Range("B15").Value = 1.401082665866568
$ws.Range("C15").Value = 0.3008255892684986
$ws.Range("D15").Value = 0.0336128826437232
$ws.Range("E15").Value = 0.1144229204760743
$ws.Range("F15").Value = 0.8661278103283507
$ws.Range("L15").Value = 0.2028822879688903
$ws.Range("M15").Value = 0.2865805410885116
$ws.Range("N15").Value = 1.318597679860268
$ws.Range("O15").Value = 2.928836463947505
$ws.Range("B16").Value = 1.334403988727388
$ws.Range("C16").Value = 0.2930510194658211
$ws.Range("D16").Value = 0.03273479691502246
$ws.Range("E16").Value = 0.1150315597160443
$ws.Range("F16").Value = 0.8590400742649393
$ws.Range("L16").Value = 0.2003856583486794
$ws.Range("M16").Value = 0.2758621278752074
$ws.Range("N16").Value = 1.327292917792221
$ws.Range("O16").Value = 2.916624646883719
$ws.Range("B17").Value = 1.293550397174499
$ws.Range("C17").Value = 0.2882664864480375
$ws.Range("D17").Value = 0.03219437568794348
$ws.Range("E17").Value = 0.1154149308538646
$ws.Range("F17").Value = 0.8548432971164885
$ws.Range("L17").Value = 0.1988784412296667
$ws.Range("M17").Value = 0.2693067203457957
$ws.Range("N17").Value = 1.332770071880859
$ws.Range("O17").Value = 2.909643114798371
$ws.Range("B18").Value = 1.270071864072975
$ws.Range("C18").Value = 0.2855090555086974
$ws.Range("D18").Value = 0.03188290485468315
$ws.Range("E18").Value = 0.1156391099730074
$ws.Range("F18").Value = 0.8524849927366915
$ws.Range("L18").Value = 0.1980205072332097
$ws.Range("M18").Value = 0.2655436339572717
$ws.Range("N18").Value = 1.335972893777637
$ws.Range("O18").Value = 2.905814794237386
$ws.Range("B19").Value = 1.262125801635136
$ws.Range("C19").Value = 0.2845744979074993
$ws.Range("D19").Value = 0.0317773377277959
$ws.Range("E19").Value = 0.1157156444658929
$ws.Range("F19").Value = 0.851696052074729
$ws.Range("L19").Value = 0.1977315689491945
$ws.Range("M19").Value = 0.2642707963951025
$ws.Range("N19").Value = 1.337066331576068
$ws.Range("O19").Value = 2.904550739695537
$ws.Range("B20").Value = 1.297897335878531
$ws.Range("C20").Value = 0.2887763779035311
$ws.Range("D20").Value = 0.03225197029040316
$ws.Range("E20").Value = 0.1153737401375583
$ws.Range("F20").Value = 0.8552842989633973
$ws.Range("L20").Value = 0.1990379583027675
$ws.Range("M20").Value = 0.2700037891322395
$ws.Range("N20").Value = 1.332181584900418
$ws.Range("O20").Value = 2.910366925250315
$ws.Range("B21").Value = 1.418313419164519
$ws.Range("C21").Value = 0.3028281255290608
$ws.Range("D21").Value = 0.03383904294741313
$ws.Range("E21").Value = 0.1142688683651514
$ws.Range("F21").Value = 0.8680045376944179
$ws.Range("L21").Value = 0.2035343940470682
$ws.Range("M21").Value = 0.289353953445719
$ws.Range("N21").Value = 1.316396930345135
$ws.Range("O21").Value = 2.932147123433481
$ws.Range("B22").Value = 1.497176238783595
$ws.Range("C22").Value = 0.3119617681347506
$ws.Range("D22").Value = 0.03487050392575952
$ws.Range("E22").Value = 0.1135794330497495
$ws.Range("F22").Value = 0.8768136165092386
$ws.Range("L22").Value = 0.2065526718714779
$ws.Range("M22").Value = 0.302065042221912
$ws.Range("N22").Value = 1.306548548254383
$ws.Range("O22").Value = 2.948052759101927
$ws.Range("B23").Value = 1.455071080470532
$ws.Range("C23").Value = 0.3070915963349137
$ws.Range("D23").Value = 0.03432052917937511
$ws.Range("E23").Value = 0.1139444117434643
$ws.Range("F23").Value = 0.8720666524777698
$ws.Range("L23").Value = 0.2049344937795468
$ws.Range("M23").Value = 0.2952750460107012
$ws.Range("N23").Value = 1.311761999376259
$ws.Range("O23").Value = 2.939410560490671
$ws.Range("B24").Value = 1.29593205756845
$ws.Range("C24").Value = 0.2885458769583238
$ws.Range("D24").Value = 0.03222593418452391
$ws.Range("E24").Value = 0.1153923507044112
$ws.Range("F24").Value = 0.8550847523390104
$ws.Range("L24").Value = 0.1989658138887052
$ws.Range("M24").Value = 0.2696886265495309
$ws.Range("N24").Value = 1.332447471864882
$ws.Range("O24").Value = 2.910039112895078
$ws.Range("B25").Value = 1.125019938851722
$ws.Range("C25").Value = 0.268331829355418
$ws.Range("D25").Value = 0.02994237257367871
$ws.Range("E25").Value = 0.1170946956531704
$ws.Range("F25").Value = 0.8388910036447612
$ws.Range("L25").Value = 0.1928710261636013
$ws.Range("M25").Value = 0.2423735511173675
$ws.Range("N25").Value = 1.356767530887716
$ws.Range("O25").Value = 2.885514754911839
